$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 847
$ws.Range("F7").Value = 91
$ws.Range("F8").Value = 304
$ws.Range("F9").Value = 586
$ws.Range("F10").Value = 3449
$ws.Range("F12").Value = 1125
$ws.Range("F13").Value = 1562
$ws.Range("F15").Value = 869
$ws.Range("F17").Value = 1092
$ws.Range("F18").Value = 1764
$ws.Range("F21").Value = 1527
$ws.Range("F23").Value = 904
$ws.Range("F24").Value = 137
$ws.Range("F25").Value = 4159

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F19").Value = 21
$ws.Range("F20").Value = 21
$ws.Range("F39").Value = 407
$ws.Range("F47").Value = 73
$ws.Range("F48").Value = 21
$ws.Range("F49").Value = 21

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2523
$ws.Range("F6").Value = 2531
$ws.Range("F7").Value = 9565
$ws.Range("F8").Value = 150
$ws.Range("F11").Value = 379
$ws.Range("F12").Value = 2912
$ws.Range("F13").Value = 428
$ws.Range("F14").Value = 765
$ws.Range("F15").Value = 156

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2523
$ws.Range("F4").Value = 150
$ws.Range("F8").Value = 2912
$ws.Range("F10").Value = 765
$ws.Range("F13").Value = 91
$ws.Range("F14").Value = 304
$ws.Range("F15").Value = 586
$ws.Range("F18").Value = 1125
$ws.Range("F21").Value = 869
$ws.Range("F24").Value = 1092
$ws.Range("F27").Value = 21
$ws.Range("F28").Value = 21
$ws.Range("F32").Value = 1764
$ws.Range("F35").Value = 1527
$ws.Range("F41").Value = 904
$ws.Range("F43").Value = 137
$ws.Range("F45").Value = 4159
$ws.Range("F46").Value = 407
$ws.Range("F50").Value = 73
$ws.Range("F51").Value = 21
$ws.Range("F52").Value = 21
